$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column K ("Ship To Customer Name") rows 2-14: tvivqt006889 -> htudix371430
$ws.Range("K2:K14").Value = "htudix371430"

# Column AX ("Previous Doc") rows 2-14: 8728482051 -> 7038567858
$ws.Range("AX2:AX14").Value = "7038567858"
